$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (shared string) into a cell without ever
# letting Excel's automatic number/percentage inference turn it into a
# numeric value (and without touching styles.xml): build the literal via
# a text formula in a scratch cell, then paste-special VALUES ONLY onto
# the destination, then clear the scratch cell.
$scratch = $ws.Range("Z1")
function Set-TextValue([string]$addr, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.ClearContents()
}

# KPI block (rows 2-7)
Set-TextValue "D2" "11"
Set-TextValue "D3" "8"
Set-TextValue "D4" "13"
Set-TextValue "D5" "16"
Set-TextValue "D6" "3"
Set-TextValue "E6" "13"
Set-TextValue "F6" "23.1%"
Set-TextValue "E7" "11"

# MATCH row (8)
Set-TextValue "D8" "2"
Set-TextValue "E8" "6"
Set-TextValue "F8" "33.3%"
Set-TextValue "G8" "4_BR,NULL_GR"

# NO MATCH row (9)
Set-TextValue "E9" "6"
Set-TextValue "F9" "66.7%"

# GAP ref row (10)
Set-TextValue "E10" "11"
Set-TextValue "F10" "45.5%"

# GAP new row (11)
Set-TextValue "E11" "8"
Set-TextValue "F11" "25.0%"

# DUPS both row (12)
Set-TextValue "D12" "1"
Set-TextValue "E12" "6"
Set-TextValue "F12" "16.7%"
Set-TextValue "G12" "NULL_GR"

# DUPS ref row (13)
Set-TextValue "E13" "11"
Set-TextValue "F13" "9.1%"

# DUPS new row (14)
Set-TextValue "E14" "8"
Set-TextValue "F14" "25.0%"
